# Auto-generated Excel COM-interop script
# Applies the KHL probabilities "human tour" publish update:
#   - Replaces match data in rows 2-4 of "Summary" and "Cards_telegram"
#   - Deletes the now-obsolete row 5 (Северсталь - Лада) from both sheets

$wb = $excel.ActiveWorkbook
$wsSummary = $wb.Worksheets.Item("Summary")
$wsCards = $wb.Worksheets.Item("Cards_telegram")

# --- Summary row 2: Локомотив – Сибирь ---
$wsSummary.Cells.Item(2, 1).Value = 1369
$wsSummary.Cells.Item(2, 2).Value = 45985.79166666666
$wsSummary.Cells.Item(2, 3).Value = 'Локомотив'
$wsSummary.Cells.Item(2, 4).Value = 'Сибирь'
$wsSummary.Cells.Item(2, 5).Value = 'Локомотив – Сибирь'
$wsSummary.Cells.Item(2, 6).Value = 897811
$wsSummary.Cells.Item(2, 7).Value = 'https://text.khl.ru/text/897811.html'
$wsSummary.Cells.Item(2, 8).Value = 3.816985
$wsSummary.Cells.Item(2, 9).Value = 1.117647
$wsSummary.Cells.Item(2, 10).Value = 4.934632
$wsSummary.Cells.Item(2, 11).Value = 35.165054
$wsSummary.Cells.Item(2, 12).Value = 20.258009
$wsSummary.Cells.Item(2, 13).Value = 55.423063
$wsSummary.Cells.Item(2, 14).Value = 0.79889
$wsSummary.Cells.Item(2, 15).Value = 0.092067
$wsSummary.Cells.Item(2, 16).Value = 0.099583
$wsSummary.Cells.Item(2, 17).Value = 1.251736784788895
$wsSummary.Cells.Item(2, 18).Value = 10.86165509900399
$wsSummary.Cells.Item(2, 19).Value = 10.04187461715353
$wsSummary.Cells.Item(2, 20).Value = 79.889
$wsSummary.Cells.Item(2, 21).Value = 9.2067
$wsSummary.Cells.Item(2, 22).Value = 9.958300000000001
$wsSummary.Cells.Item(2, 23).Value = 0.199104
$wsSummary.Cells.Item(2, 24).Value = 0.791435
$wsSummary.Cells.Item(2, 25).Value = 1.263527642826006
$wsSummary.Cells.Item(2, 26).Value = 0.336547
$wsSummary.Cells.Item(2, 27).Value = 0.653992
$wsSummary.Cells.Item(2, 28).Value = 1.529070691996232
$wsSummary.Cells.Item(2, 29).Value = 0.490706
$wsSummary.Cells.Item(2, 30).Value = 0.499833
$wsSummary.Cells.Item(2, 31).Value = 2.000668223186544
$wsSummary.Cells.Item(2, 32).Value = 0.949509
$wsSummary.Cells.Item(2, 33).Value = 0.050491
$wsSummary.Cells.Item(2, 34).Value = 1.053175904599114
$wsSummary.Cells.Item(2, 35).Value = 0.850887
$wsSummary.Cells.Item(2, 36).Value = 0.149113
$wsSummary.Cells.Item(2, 37).Value = 1.175244186360821
$wsSummary.Cells.Item(2, 38).Value = 0.593374
$wsSummary.Cells.Item(2, 39).Value = 0.406626
$wsSummary.Cells.Item(2, 40).Value = 1.685277750626081
$wsSummary.Cells.Item(2, 41).Value = 0.322704
$wsSummary.Cells.Item(2, 42).Value = 0.677296
$wsSummary.Cells.Item(2, 43).Value = 3.098815013138976
$wsSummary.Cells.Item(2, 44).Value = 0.945669
$wsSummary.Cells.Item(2, 45).Value = 1.057452449006999
$wsSummary.Cells.Item(2, 46).Value = 0.321248
$wsSummary.Cells.Item(2, 47).Value = 3.112859846598267

# --- Summary row 3: Динамо М – Амур ---
$wsSummary.Cells.Item(3, 1).Value = 1369
$wsSummary.Cells.Item(3, 2).Value = 45985.8125
$wsSummary.Cells.Item(3, 3).Value = 'Динамо М'
$wsSummary.Cells.Item(3, 4).Value = 'Амур'
$wsSummary.Cells.Item(3, 5).Value = 'Динамо М – Амур'
$wsSummary.Cells.Item(3, 6).Value = 897809
$wsSummary.Cells.Item(3, 7).Value = 'https://text.khl.ru/text/897809.html'
$wsSummary.Cells.Item(3, 8).Value = 1.983158
$wsSummary.Cells.Item(3, 9).Value = 3.411544
$wsSummary.Cells.Item(3, 10).Value = 5.394701
$wsSummary.Cells.Item(3, 11).Value = 29.463809
$wsSummary.Cells.Item(3, 12).Value = 33.087107
$wsSummary.Cells.Item(3, 13).Value = 62.550916
$wsSummary.Cells.Item(3, 14).Value = 0.396516
$wsSummary.Cells.Item(3, 15).Value = 0.161561
$wsSummary.Cells.Item(3, 16).Value = 0.44098
$wsSummary.Cells.Item(3, 17).Value = 2.521966326705606
$wsSummary.Cells.Item(3, 18).Value = 6.189612592147857
$wsSummary.Cells.Item(3, 19).Value = 2.267676538618531
$wsSummary.Cells.Item(3, 20).Value = 39.65159999999999
$wsSummary.Cells.Item(3, 21).Value = 16.1561
$wsSummary.Cells.Item(3, 22).Value = 44.098
$wsSummary.Cells.Item(3, 23).Value = 0.241218
$wsSummary.Cells.Item(3, 24).Value = 0.757839
$wsSummary.Cells.Item(3, 25).Value = 1.31954148572454
$wsSummary.Cells.Item(3, 26).Value = 0.391581
$wsSummary.Cells.Item(3, 27).Value = 0.607476
$wsSummary.Cells.Item(3, 28).Value = 1.646155568285825
$wsSummary.Cells.Item(3, 29).Value = 0.550654
$wsSummary.Cells.Item(3, 30).Value = 0.448403
$wsSummary.Cells.Item(3, 31).Value = 2.230136729682897
$wsSummary.Cells.Item(3, 32).Value = 0.815552
$wsSummary.Cells.Item(3, 33).Value = 0.184448
$wsSummary.Cells.Item(3, 34).Value = 1.226163383818567
$wsSummary.Cells.Item(3, 35).Value = 0.599231
$wsSummary.Cells.Item(3, 36).Value = 0.400769
$wsSummary.Cells.Item(3, 37).Value = 1.668805519073613
$wsSummary.Cells.Item(3, 38).Value = 0.834675
$wsSummary.Cells.Item(3, 39).Value = 0.165325
$wsSummary.Cells.Item(3, 40).Value = 1.198071105520113
$wsSummary.Cells.Item(3, 41).Value = 0.629565
$wsSummary.Cells.Item(3, 42).Value = 0.370435
$wsSummary.Cells.Item(3, 43).Value = 1.588398338535338
$wsSummary.Cells.Item(3, 44).Value = 0.709696
$wsSummary.Cells.Item(3, 45).Value = 1.409054017494815
$wsSummary.Cells.Item(3, 46).Value = 0.747436
$wsSummary.Cells.Item(3, 47).Value = 1.337907192053902

# --- Summary row 4: ЦСКА – СКА ---
$wsSummary.Cells.Item(4, 1).Value = 1369
$wsSummary.Cells.Item(4, 2).Value = 45985.8125
$wsSummary.Cells.Item(4, 3).Value = 'ЦСКА'
$wsSummary.Cells.Item(4, 4).Value = 'СКА'
$wsSummary.Cells.Item(4, 5).Value = 'ЦСКА – СКА'
$wsSummary.Cells.Item(4, 6).Value = 897810
$wsSummary.Cells.Item(4, 7).Value = 'https://text.khl.ru/text/897810.html'
$wsSummary.Cells.Item(4, 8).Value = 4.615385
$wsSummary.Cells.Item(4, 9).Value = 2.109924
$wsSummary.Cells.Item(4, 10).Value = 6.725309
$wsSummary.Cells.Item(4, 11).Value = 34.357697
$wsSummary.Cells.Item(4, 12).Value = 28.798769
$wsSummary.Cells.Item(4, 13).Value = 63.156466
$wsSummary.Cells.Item(4, 14).Value = 0.676504
$wsSummary.Cells.Item(4, 15).Value = 0.151808
$wsSummary.Cells.Item(4, 16).Value = 0.171278
$wsSummary.Cells.Item(4, 17).Value = 1.478187859938744
$wsSummary.Cells.Item(4, 18).Value = 6.587268128161889
$wsSummary.Cells.Item(4, 19).Value = 5.838461448639054
$wsSummary.Cells.Item(4, 20).Value = 67.6504
$wsSummary.Cells.Item(4, 21).Value = 15.1808
$wsSummary.Cells.Item(4, 22).Value = 17.1278
$wsSummary.Cells.Item(4, 23).Value = 0.487857
$wsSummary.Cells.Item(4, 24).Value = 0.511733
$wsSummary.Cells.Item(4, 25).Value = 1.95414405559149
$wsSummary.Cells.Item(4, 26).Value = 0.662067
$wsSummary.Cells.Item(4, 27).Value = 0.337523
$wsSummary.Cells.Item(4, 28).Value = 2.962761056283572
$wsSummary.Cells.Item(4, 29).Value = 0.799595
$wsSummary.Cells.Item(4, 30).Value = 0.199996
$wsSummary.Cells.Item(4, 31).Value = 5.00010000200004
$wsSummary.Cells.Item(4, 32).Value = 0.818617
$wsSummary.Cells.Item(4, 33).Value = 0.181383
$wsSummary.Cells.Item(4, 34).Value = 1.221572481392397
$wsSummary.Cells.Item(4, 35).Value = 0.604001
$wsSummary.Cells.Item(4, 36).Value = 0.395999
$wsSummary.Cells.Item(4, 37).Value = 1.65562639796954
$wsSummary.Cells.Item(4, 38).Value = 0.479147
$wsSummary.Cells.Item(4, 39).Value = 0.520853
$wsSummary.Cells.Item(4, 40).Value = 2.087042181209525
$wsSummary.Cells.Item(4, 41).Value = 0.219917
$wsSummary.Cells.Item(4, 42).Value = 0.780083
$wsSummary.Cells.Item(4, 43).Value = 4.547170068707739
$wsSummary.Cells.Item(4, 44).Value = 0.924297
$wsSummary.Cells.Item(4, 45).Value = 1.081903327610065
$wsSummary.Cells.Item(4, 46).Value = 0.509023
$wsSummary.Cells.Item(4, 47).Value = 1.964547770925872

# Remove the 5th match (Северсталь - Лада) - no longer part of this publish batch
$wsSummary.Rows.Item(5).EntireRow.Delete()

# --- Cards_telegram row 2: Локомотив – Сибирь ---
$wsCards.Cells.Item(2, 1).Value = 45985.79166666666
$wsCards.Cells.Item(2, 2).Value = 'Локомотив – Сибирь'
$card2 = @"
КХЛ • Регулярный чемпионат • 24.11.2025
Локомотив – Сибирь
Ожидания модели (60’):
• Голы: λ_total ≈ 6.73 (4.73 : 2.00)
• Броски: SOG λ ≈ 55 (35 : 20)
Исход (60’), честные кф:
• П1: 79.9%  (Kмод 1.25)
• Х:  9.2%  (Kмод 10.86)
• П2: 10.0%  (Kмод 10.04)
Тоталы голов:
• ТМ 4.5: 19.9%  (Kмод 5.02)
• ТБ 4.5: 79.1%  (Kмод 1.26)
• ТМ 5.5: 33.7%  (Kмод 2.97)
• ТБ 5.5: 65.4%  (Kмод 1.53)
• ТМ 6.5: 49.1%  (Kмод 2.04)
• ТБ 6.5: 50.0%  (Kмод 2.00)
Индивидуальные тоталы:
• Локомотив ИТБ 1.5: 95.0% (Kмод 1.05)
• Локомотив ИТБ 2.5: 85.1% (Kмод 1.18)
• Сибирь ИТБ 1.5: 59.3% (Kмод 1.69)
• Сибирь ИТБ 2.5: 32.3% (Kмод 3.10)
Фора +1.5:
• Локомотив +1.5: 94.6% (Kмод 1.06)
• Сибирь +1.5: 32.1% (Kмод 3.11)
"@
$wsCards.Cells.Item(2, 3).Value = $card2

# --- Cards_telegram row 3: Динамо М – Амур ---
$wsCards.Cells.Item(3, 1).Value = 45985.8125
$wsCards.Cells.Item(3, 2).Value = 'Динамо М – Амур'
$card3 = @"
КХЛ • Регулярный чемпионат • 24.11.2025
Динамо М – Амур
Ожидания модели (60’):
• Голы: λ_total ≈ 6.35 (3.10 : 3.25)
• Броски: SOG λ ≈ 63 (29 : 33)
Исход (60’), честные кф:
• П1: 39.7%  (Kмод 2.52)
• Х:  16.2%  (Kмод 6.19)
• П2: 44.1%  (Kмод 2.27)
Тоталы голов:
• ТМ 4.5: 24.1%  (Kмод 4.15)
• ТБ 4.5: 75.8%  (Kмод 1.32)
• ТМ 5.5: 39.2%  (Kмод 2.55)
• ТБ 5.5: 60.7%  (Kмод 1.65)
• ТМ 6.5: 55.1%  (Kмод 1.82)
• ТБ 6.5: 44.8%  (Kмод 2.23)
Индивидуальные тоталы:
• Динамо М ИТБ 1.5: 81.6% (Kмод 1.23)
• Динамо М ИТБ 2.5: 59.9% (Kмод 1.67)
• Амур ИТБ 1.5: 83.5% (Kмод 1.20)
• Амур ИТБ 2.5: 63.0% (Kмод 1.59)
Фора +1.5:
• Динамо М +1.5: 71.0% (Kмод 1.41)
• Амур +1.5: 74.7% (Kмод 1.34)
"@
$wsCards.Cells.Item(3, 3).Value = $card3

# --- Cards_telegram row 4: ЦСКА – СКА ---
$wsCards.Cells.Item(4, 1).Value = 45985.8125
$wsCards.Cells.Item(4, 2).Value = 'ЦСКА – СКА'
$card4 = @"
КХЛ • Регулярный чемпионат • 24.11.2025
ЦСКА – СКА
Ожидания модели (60’):
• Голы: λ_total ≈ 4.74 (3.12 : 1.61)
• Броски: SOG λ ≈ 63 (34 : 29)
Исход (60’), честные кф:
• П1: 67.7%  (Kмод 1.48)
• Х:  15.2%  (Kмод 6.59)
• П2: 17.1%  (Kмод 5.84)
Тоталы голов:
• ТМ 4.5: 48.8%  (Kмод 2.05)
• ТБ 4.5: 51.2%  (Kмод 1.95)
• ТМ 5.5: 66.2%  (Kмод 1.51)
• ТБ 5.5: 33.8%  (Kмод 2.96)
• ТМ 6.5: 80.0%  (Kмод 1.25)
• ТБ 6.5: 20.0%  (Kмод 5.00)
Индивидуальные тоталы:
• ЦСКА ИТБ 1.5: 81.9% (Kмод 1.22)
• ЦСКА ИТБ 2.5: 60.4% (Kмод 1.66)
• СКА ИТБ 1.5: 47.9% (Kмод 2.09)
• СКА ИТБ 2.5: 22.0% (Kмод 4.55)
Фора +1.5:
• ЦСКА +1.5: 92.4% (Kмод 1.08)
• СКА +1.5: 50.9% (Kмод 1.96)
"@
$wsCards.Cells.Item(4, 3).Value = $card4

# Remove the 5th card (Северсталь - Лада) - no longer part of this publish batch
$wsCards.Rows.Item(5).EntireRow.Delete()

